# CRD_DTR_Flow.pptx edit
#
# On slide 3, the "Rectangle 38" label shape reads "DTR SMART Client" on a
# single line. It is changed to wrap across two lines/paragraphs:
#   "SMART "
#   "DTR Client"
# Inserting a carriage return inside the TextRange splits the run into two
# paragraphs, and PowerPoint carries the existing run formatting (font,
# size, color, etc.) over onto the new paragraph automatically.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$target = $null
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Rectangle 38") {
        $target = $shp
        break
    }
}

if ($target -eq $null) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "DTR SMART Client") {
            $target = $shp
            break
        }
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "SMART `rDTR Client"
}
